# Update Name of Algo
# Apply the updated algorithm result values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = 12.914
$ws.Range("D9").Value  = -7.998
$ws.Range("E11").Value = 12.914
$ws.Range("D18").Value = -8.144
$ws.Range("D20").Value = -8.044
